# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.878.06"
$ws.Range("E2").Value = "'  +1.78%  "

# Row 3
$ws.Range("D3").Value = "'1.867.99"
$ws.Range("E3").Value = "'  +1.67%  "

# Row 4
$ws.Range("E4").Value = "'  +0.87%  "

# Row 5
$ws.Range("D5").Value = "'323.29"
$ws.Range("E5").Value = "'  +1.55%  "

# Row 6
$ws.Range("E6").Value = "'  +0.40%  "

# Row 7
$ws.Range("D7").Value = "'0.4432"
$ws.Range("E7").Value = "'  +1.75%  "

# Row 8
$ws.Range("D8").Value = "'0.3822"
$ws.Range("E8").Value = "'  +2.62%  "

# Row 9
$ws.Range("D9").Value = "'0.07481"
$ws.Range("E9").Value = "'  +2.06%  "

# Row 10
$ws.Range("D10").Value = "'0.8909"
$ws.Range("E10").Value = "'  +2.18%  "

# Row 11
$ws.Range("D11").Value = "'21.82"
$ws.Range("E11").Value = "'  +2.36%  "

# Row 12
$ws.Range("D12").Value = "'1.897.57"
$ws.Range("E12").Value = "'  -2.68%  "

# Row 13
$ws.Range("D13").Value = "'5.575"
$ws.Range("E13").Value = "'  +2.02%  "

# Row 14
$ws.Range("D14").Value = "'6.787"
$ws.Range("E14").Value = "'  +1.70%  "

# Row 15
$ws.Range("D15").Value = "'0.07203"
$ws.Range("E15").Value = "'  +1.10%  "

# Row 16
$ws.Range("D16").Value = "'84.82"
$ws.Range("E16").Value = "'  +3.28%  "

# Row 17
$ws.Range("D17").Value = "'1.036"
$ws.Range("E17").Value = "'  +0.60%  "

# Row 18
$ws.Range("D18").Value = "'0.000009149"
$ws.Range("E18").Value = "'  +1.97%  "

# Row 19
$ws.Range("D19").Value = "'1.029"
$ws.Range("E19").Value = "'  +0.48%  "

# Row 20
$ws.Range("D20").Value = "'15.60"
$ws.Range("E20").Value = "'  +1.46%  "

# Row 21
$ws.Range("D21").Value = "'27.892.54"
$ws.Range("E21").Value = "'  +1.70%  "

# Row 22
$ws.Range("D22").Value = "'5.329"
$ws.Range("E22").Value = "'  +1.62%  "

# Row 23
$ws.Range("D23").Value = "'11.34"
$ws.Range("E23").Value = "'  +1.49%  "

# Row 24
$ws.Range("D24").Value = "'2.103.91"
$ws.Range("E24").Value = "'  -1.91%  "

# Row 25
$ws.Range("D25").Value = "'2.021"
$ws.Range("E25").Value = "'  +6.82%  "

# Row 26
$ws.Range("D26").Value = "'158.34"
$ws.Range("E26").Value = "'  +1.03%  "

# Row 27
$ws.Range("D27").Value = "'18.91"
$ws.Range("E27").Value = "'  +2.07%  "

# Row 28
$ws.Range("D28").Value = "'5.396"
$ws.Range("E28").Value = "'  +2.97%  "

# Row 29
$ws.Range("D29").Value = "'1.987"
$ws.Range("E29").Value = "'  +3.79%  "

# Row 30
$ws.Range("D30").Value = "'118.44"
$ws.Range("E30").Value = "'  +2.74%  "

# Row 31
$ws.Range("D31").Value = "'0.09072"
$ws.Range("E31").Value = "'  +0.35%  "

# Row 32
$ws.Range("D32").Value = "'1.238"
$ws.Range("E32").Value = "'  +3.45%  "

# Row 33
$ws.Range("D33").Value = "'0.7815"
$ws.Range("E33").Value = "'  +3.23%  "

# Row 34
$ws.Range("D34").Value = "'4.614"
$ws.Range("E34").Value = "'  +3.56%  "

# Row 35
$ws.Range("D35").Value = "'3.002"
$ws.Range("E35").Value = "'  +4.97%  "

# Row 36
$ws.Range("D36").Value = "'1.031"
$ws.Range("E36").Value = "'  +0.63%  "

# Row 37
$ws.Range("D37").Value = "'1.147"
$ws.Range("E37").Value = "'  -0.04%  "

# Row 38
$ws.Range("B38").Value = "'Hedera"
$ws.Range("C38").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05372"
$ws.Range("E38").Value = "'  +2.63%  "

# Row 39
$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01989"
$ws.Range("E39").Value = "'  +1.83%  "

# Row 40
$ws.Range("D40").Value = "'2.888"
$ws.Range("E40").Value = "'  +3.97%  "

# Row 41
$ws.Range("D41").Value = "'0.5225"
$ws.Range("E41").Value = "'  +1.13%  "

# Row 42
$ws.Range("D42").Value = "'0.1698"
$ws.Range("E42").Value = "'  +2.30%  "

# Row 43
$ws.Range("D43").Value = "'6.910"
$ws.Range("E43").Value = "'  +5.86%  "

# Row 44
$ws.Range("D44").Value = "'8.832"
$ws.Range("E44").Value = "'  +4.58%  "

# Row 45
$ws.Range("D45").Value = "'112.08"
$ws.Range("E45").Value = "'  +3.46%  "

# Row 46
$ws.Range("D46").Value = "'10.84"
$ws.Range("E46").Value = "'  +3.71%  "

# Row 47
$ws.Range("D47").Value = "'0.06633"
$ws.Range("E47").Value = "'  +5.55%  "

# Row 48
$ws.Range("B48").Value = "'NEARProtocol"
$ws.Range("C48").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.724"
$ws.Range("E48").Value = "'  +3.42%  "

# Row 49
$ws.Range("B49").Value = "'PaxDollar"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.031"
$ws.Range("E49").Value = "'  +0.67%  "

# Row 50
$ws.Range("D50").Value = "'0.4743"
$ws.Range("E50").Value = "'  +2.55%  "

# Row 51
$ws.Range("D51").Value = "'1.918"
$ws.Range("E51").Value = "'  +2.69%  "

